$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update B14 (tray length) -> ripples into C14, C18, C23 formulas ---
$ws.Range("B14").Value = 92

# --- Add new row 24: Phenolic tube density (g/m) = 203 ---
$ws.Range("A24").Value = "Phenolic tube density (g/m)"
$ws.Range("B24").Value = 203

# --- Move the active selection to B15, matching the saved view state ---
$ws.Range("B15").Select()
